$d = $word.ActiveDocument

# 1 & 7. Title / heading text and the later bold restatement of it (both occurrences identical)
$d.Content.Find.Execute("Play Free Dazzling Diamonds Slot Game | Review & Rating", $true, $false, $false, $false, $false, $true, 1, $false, "Play Dazzling Diamonds free - Review", 2)

# 2. "What we like" bullet 2
$d.Content.Find.Execute("Crisp graphics with bright and colorful symbols", $true, $false, $false, $false, $false, $true, 1, $false, "Impressive 100,000 euro jackpot", 2)

# 3. "What we like" bullet 3
$d.Content.Find.Execute("Offers a jackpot of 100,000 euros", $true, $false, $false, $false, $false, $true, 1, $false, "Crisp design without distractions", 2)

# 4. "What we like" bullet 4
$d.Content.Find.Execute("Runs smoothly on desktop and mobile devices", $true, $false, $false, $false, $false, $true, 1, $false, "Accessible on desktop and mobile devices", 2)

# 5. "What we don't like" bullet 1
$d.Content.Find.Execute("Does not have a Wild symbol", $true, $false, $false, $false, $false, $true, 1, $false, "Lack of a Wild symbol", 2)

# 6. "What we don't like" bullet 2
$d.Content.Find.Execute("Limited range of bets", $true, $false, $false, $false, $false, $true, 1, $false, "Basic sound effects", 2)

# 8. Meta description (italic run)
$d.Content.Find.Execute("Play the free Dazzling Diamonds slot game with a jackpot of 100,000 euros. Read our review on features, graphics, payout potential, accessibility, and more.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Dazzling Diamonds slot game and play for free. Enjoy simple gameplay and an impressive 100,000 euro jackpot.", 2)
